# Case-closure template: the "services_provided" block-helper delimiters
# ({-w:p services_provided} ... {/services_provided}) were wrongly
# italicised along with the template expressions they wrap. Strip the
# italics from just the delimiter/helper-name text while leaving the
# {service} / {description} placeholders themselves untouched, so the
# merge field renders correctly when the form is downloaded.

$d = $word.ActiveDocument

# --- Opening delimiter: "{-w:p services_provided}" -------------------
$opening = $d.Content
$opening.Find.ClearFormatting()
$foundOpening = $opening.Find.Execute(
    "{-w:p services_provided}", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $foundOpening) {
    throw "Could not find the opening services_provided delimiter"
}
$opening.Font.Italic = $false

# --- Closing delimiter: "{/services_provided}" ------------------------
$closing = $d.Content
$closing.Find.ClearFormatting()
$foundClosing = $closing.Find.Execute(
    "{/services_provided}", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $foundClosing) {
    throw "Could not find the closing services_provided delimiter"
}
$closing.Font.Italic = $false
